# This script inserts a new asset class, "returns Private Equity USD Unhedged",
# into the correlation matrix on Sheet1. It is placed alphabetically between
# "returns Nature Capital - USD Hedged" (col G / row 7) and
# "returns UK Property Direct - USD Unhedged" (previously col H / row 8),
# so a new column H and a new row 8 are inserted, pushing everything else
# one column/row to the right/down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new column (H) and new row (8) ----------------------------
$ws.Columns("H:H").Insert()
$ws.Rows("8:8").Insert()

# --- Carry over formatting for the new header cell and new row label ------
# Header row style (bold, centered, bordered) lives on row 1; copy it from
# the neighbouring header cell G1 into the freshly inserted H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Row-label style lives in column A; copy it from the neighbouring label
# cell A7 into the freshly inserted A8.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Header / row label text ------------------------------------------------
$ws.Range("H1").Value = "returns Private Equity USD Unhedged"
$ws.Range("A8").Value = "returns Private Equity USD Unhedged"

# --- New column H correlation values (rows 2-10) ---------------------------
$ws.Range("H2").Value  = 0.5444423922450708
$ws.Range("H3").Value  = 0.9151257149970743
$ws.Range("H4").Value  = 0.8127456037382913
$ws.Range("H5").Value  = -0.1280087661092652
$ws.Range("H6").Value  = 0.7325619217938056
$ws.Range("H7").Value  = 0.354584949574704
$ws.Range("H8").Value  = 1
$ws.Range("H9").Value  = 0.8552535606935843
$ws.Range("H10").Value = 0.875528909419702

# --- New row 8 correlation values (cols B-G, I, J) --------------------------
$ws.Range("B8").Value = 0.5444423922450708
$ws.Range("C8").Value = 0.9151257149970743
$ws.Range("D8").Value = 0.8127456037382913
$ws.Range("E8").Value = -0.1280087661092652
$ws.Range("F8").Value = 0.7325619217938056
$ws.Range("G8").Value = 0.354584949574704
$ws.Range("I8").Value = 0.8552535606935843
$ws.Range("J8").Value = 0.875528909419702
